$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction in SA algorithm and 746 logs:
# Fitness values (column C) for Generation 0-87 (rows 2-89) are corrected to 7590
# Fitness values (column C) for Generation 88-115 (rows 90-117) are corrected to 7573
# Rows 118+ already hold the corrected value (7573) and are left untouched.

$ws.Range("C2:C89").Value = 7590
$ws.Range("C90:C117").Value = 7573
